$d = $word.ActiveDocument

# Locate the run containing "Selecionar Tarefa" (table cell: "Nome do requisito: Selecionar Tarefa")
$rng = $d.Content
$found = $rng.Find.Execute("Selecionar Tarefa")

if ($found) {
    # Shrink the run's text down to "Escolher" (keeps it as its own run, same formatting
    # as before, distinct from the preceding "Nome do requisito: " run).
    $rng.Text = "Escolher"

    # Temporarily perturb formatting so the upcoming insertion does not get silently
    # re-merged with the (identically formatted) "Escolher" run once we restore it.
    $rng.Font.Color = 255

    # Insert " Tarefa" right after "Escolher" as a brand new run.
    $rng.Collapse(0)
    $rng.InsertAfter(" Tarefa")

    # $rng now spans the newly inserted " Tarefa" text -- restore its real color (black).
    $rng.Font.Color = 0

    # Restore "Escolher" back to black as well, as its own separate run.
    $rngEscolher = $d.Content
    $rngEscolher.Find.Execute("Escolher") | Out-Null
    $rngEscolher.Font.Color = 0
}
